# Mise à jour de certains champs de Modules et de Professeurs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates (C1/D1 swap content + wording changes):
#   C1: "Enseignant"      -> "Chef  Module"
#   D1: "Nombre d'heures" -> "Composants"
$ws.Range("C1").Value2 = "Chef  Module"
$ws.Range("D1").Value2 = "Composants"

# Column widths for C and D (engine rounds ColumnWidth to the nearest
# achievable pixel grid step, so we pick the input that lands closest to
# the target stored width: C -> 35 exactly, D -> 24.5703125 (~24.5)).
$ws.Columns("C").ColumnWidth = 34.166666666666664
$ws.Columns("D").ColumnWidth = 23.666666666666668

# Update the current selection to E8
$null = $ws.Range("E8").Select()
